$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New Cypher query text for the "ParticipantsTab" row (B2), replacing the old
# primary-diagnosis participant query with the corrected/fixed version.
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in ['Acinar cell carcinoma']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$ws.Range("B2").Value = $newParticipantQuery

# The longer query text wraps onto more lines, so the row needs to grow to
# fit it (mirrors Excel's automatic row-height recalculation for wrapped
# text rows).
$ws.Rows.Item(2).RowHeight = 330.75

# Update the selected cell shown when the sheet is reopened.
$ws.Range("B2").Select()
